$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scalings")

# Warrior's max_hp / phy_atk / mag_atk scaling values were retuned
$ws.Range("B2").Value = "90 * Lvl"
$ws.Range("C2").Value = "8 * Lvl"
$ws.Range("D2").Value = "8 * Lvl"

# A2 ("Warrior") no longer carries the small-font override style
$ws.Range("A2").ClearFormats()

# Widen column A so the class names aren't clipped
$ws.Columns("A").ColumnWidth = 19.28

# Move the active selection
[void]$ws.Range("C10").Select()
